# Case_5_124 res_line/pl_mw.xlsx update: "case with 380 kV done"
# Updates columns B, C, D, E, F, H, M for data rows 2-25 (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.947645898192377
$ws.Cells.Item(2, 3).Value = 0.1719206942726998
$ws.Cells.Item(2, 4).Value = 0.1443690170293905
$ws.Cells.Item(2, 5).Value = 0.0911239541133213
$ws.Cells.Item(2, 6).Value = 3.037772572775197
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 13).Value = 0.3725926983130847

$ws.Cells.Item(3, 2).Value = 0.8641610323122677
$ws.Cells.Item(3, 3).Value = 0.150468738092485
$ws.Cells.Item(3, 4).Value = 0.1336508942512609
$ws.Cells.Item(3, 5).Value = 0.08435741998362545
$ws.Cells.Item(3, 6).Value = 2.820701146081205
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 13).Value = 0.3392264900062472

$ws.Cells.Item(4, 2).Value = 0.8137861859623001
$ws.Cells.Item(4, 3).Value = 0.1374226980947242
$ws.Cells.Item(4, 4).Value = 0.127071828517586
$ws.Cells.Item(4, 5).Value = 0.08024515227909035
$ws.Cells.Item(4, 6).Value = 2.688233000857025
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 13).Value = 0.3190326487632831

$ws.Cells.Item(5, 2).Value = 0.7934772783993367
$ws.Cells.Item(5, 3).Value = 0.1321368024237302
$ws.Cells.Item(5, 4).Value = 0.1243906508097012
$ws.Cells.Item(5, 5).Value = 0.07857969643593066
$ws.Cells.Item(5, 6).Value = 2.634446044827456
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 13).Value = 0.3108757524015573

$ws.Cells.Item(6, 2).Value = 0.7901181458480835
$ws.Cells.Item(6, 3).Value = 0.131260892179256
$ws.Cells.Item(6, 4).Value = 0.1239454154905246
$ws.Cells.Item(6, 5).Value = 0.07830376266251449
$ws.Cells.Item(6, 6).Value = 2.625526243975798
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 13).Value = 0.3095256281440015

$ws.Cells.Item(7, 2).Value = 0.8135114095485676
$ws.Cells.Item(7, 3).Value = 0.1373512887153652
$ws.Cells.Item(7, 4).Value = 0.1270356706803284
$ws.Cells.Item(7, 5).Value = 0.08022264999092954
$ws.Cells.Item(7, 6).Value = 2.687506835130478
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 13).Value = 0.318922351340909

$ws.Cells.Item(8, 2).Value = 0.9186745346849534
$ws.Cells.Item(8, 3).Value = 0.1644973286930451
$ws.Cells.Item(8, 4).Value = 0.1406726273254151
$ws.Cells.Item(8, 5).Value = 0.08878184398055566
$ws.Cells.Item(8, 6).Value = 2.962751197177766
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 13).Value = 0.3610262927193943

$ws.Cells.Item(9, 2).Value = 1.132087899232772
$ws.Cells.Item(9, 3).Value = 0.2187797939818665
$ws.Cells.Item(9, 4).Value = 0.1674595635135319
$ws.Cells.Item(9, 5).Value = 0.1059191192271101
$ws.Cells.Item(9, 6).Value = 3.509445478715264
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 13).Value = 0.4459890776067965

$ws.Cells.Item(10, 2).Value = 1.293515007262101
$ws.Cells.Item(10, 3).Value = 0.2593791426511984
$ws.Cells.Item(10, 4).Value = 0.1872102887500091
$ws.Cells.Item(10, 5).Value = 0.1187487321034979
$ws.Cells.Item(10, 6).Value = 3.916036404187508
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 13).Value = 0.5099789378127184

$ws.Cells.Item(11, 2).Value = 1.368014845915184
$ws.Cells.Item(11, 3).Value = 0.2780223861166178
$ws.Cells.Item(11, 4).Value = 0.1962201569249089
$ws.Cells.Item(11, 5).Value = 0.1246424250083251
$ws.Cells.Item(11, 6).Value = 4.102233106722792
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 13).Value = 0.539453944972621

$ws.Cells.Item(12, 2).Value = 1.396383658375896
$ws.Cells.Item(12, 3).Value = 0.2851085723081894
$ws.Cells.Item(12, 4).Value = 0.199636325662766
$ws.Cells.Item(12, 5).Value = 0.1268828873086605
$ws.Cells.Item(12, 6).Value = 4.172930926557683
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 13).Value = 0.5506698174898901

$ws.Cells.Item(13, 2).Value = 1.390266863196132
$ws.Cells.Item(13, 3).Value = 0.2835812401042119
$ws.Cells.Item(13, 4).Value = 0.198900388762894
$ws.Cells.Item(13, 5).Value = 0.1263999727389233
$ws.Cells.Item(13, 6).Value = 4.157696289548142
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 13).Value = 0.5482518343006433

$ws.Cells.Item(14, 2).Value = 1.370345588116095
$ws.Cells.Item(14, 3).Value = 0.278604834273807
$ws.Cells.Item(14, 4).Value = 0.1965011163647432
$ws.Cells.Item(14, 5).Value = 0.1248265736014886
$ws.Cells.Item(14, 6).Value = 4.10804560182288
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 13).Value = 0.5403755826418006

$ws.Cells.Item(15, 2).Value = 1.358163845118554
$ws.Cells.Item(15, 3).Value = 0.2755601190016819
$ws.Cells.Item(15, 4).Value = 0.1950320774745649
$ws.Cells.Item(15, 5).Value = 0.1238639591919863
$ws.Cells.Item(15, 6).Value = 4.077658098306074
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 13).Value = 0.5355582784261657

$ws.Cells.Item(16, 2).Value = 1.2886680607117
$ws.Cells.Item(16, 3).Value = 0.2581643839288574
$ws.Cells.Item(16, 4).Value = 0.1866220393957292
$ws.Cells.Item(16, 5).Value = 0.1183647543156567
$ws.Cells.Item(16, 6).Value = 3.903893859664663
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 13).Value = 0.5080601833441705

$ws.Cells.Item(17, 2).Value = 1.246310498272237
$ws.Cells.Item(17, 3).Value = 0.2475382572563944
$ws.Cells.Item(17, 4).Value = 0.1814696997102487
$ws.Cells.Item(17, 5).Value = 0.1150061561221705
$ws.Cells.Item(17, 6).Value = 3.797619838831281
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 13).Value = 0.4912858816138481

$ws.Cells.Item(18, 2).Value = 1.22204761993595
$ws.Cells.Item(18, 3).Value = 0.2414427244910939
$ws.Cells.Item(18, 4).Value = 0.1785085358093994
$ws.Cells.Item(18, 5).Value = 0.1130797618720294
$ws.Cells.Item(18, 6).Value = 3.736609521624132
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 13).Value = 0.4816720616151855

$ws.Cells.Item(19, 2).Value = 1.213849691232497
$ws.Cells.Item(19, 3).Value = 0.239381647683814
$ws.Cells.Item(19, 4).Value = 0.1775063145926623
$ws.Cells.Item(19, 5).Value = 0.1124284306039414
$ws.Cells.Item(19, 6).Value = 3.715972003021477
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 13).Value = 0.4784228286096948

$ws.Cells.Item(20, 2).Value = 1.250809146951781
$ws.Cells.Item(20, 3).Value = 0.2486677266169863
$ws.Cells.Item(20, 4).Value = 0.1820179292708985
$ws.Cells.Item(20, 5).Value = 0.115363124757728
$ws.Cells.Item(20, 6).Value = 3.80892083470809
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 13).Value = 0.4930679707592134

$ws.Cells.Item(21, 2).Value = 1.376192646304105
$ws.Cells.Item(21, 3).Value = 0.2800657982991197
$ws.Cells.Item(21, 4).Value = 0.1972057176922988
$ws.Cells.Item(21, 5).Value = 0.1252884812071002
$ws.Cells.Item(21, 6).Value = 4.122623986738404
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 13).Value = 0.5426875401806939

$ws.Cells.Item(22, 2).Value = 1.459057083136884
$ws.Cells.Item(22, 3).Value = 0.3007407355483451
$ws.Cells.Item(22, 4).Value = 0.2071573736807863
$ws.Cells.Item(22, 5).Value = 0.1318258725214534
$ws.Cells.Item(22, 6).Value = 4.328755178232086
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 13).Value = 0.5754342946179776

$ws.Cells.Item(23, 2).Value = 1.414745303010363
$ws.Cells.Item(23, 3).Value = 0.2896915535640119
$ws.Cells.Item(23, 4).Value = 0.2018434122698523
$ws.Cells.Item(23, 5).Value = 0.1283319836148209
$ws.Cells.Item(23, 6).Value = 4.218633886059422
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 13).Value = 0.5579270984155329

$ws.Cells.Item(24, 2).Value = 1.248775030824447
$ws.Cells.Item(24, 3).Value = 0.2481570513538145
$ws.Cells.Item(24, 4).Value = 0.1817700716726165
$ws.Cells.Item(24, 5).Value = 0.1152017252236348
$ws.Cells.Item(24, 6).Value = 3.803811381126224
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 13).Value = 0.4922621951015742

$ws.Cells.Item(25, 2).Value = 1.073558865654263
$ws.Cells.Item(25, 3).Value = 0.2039746313485296
$ws.Cells.Item(25, 4).Value = 0.1602037670672019
$ws.Cells.Item(25, 5).Value = 0.101242830690488
$ws.Cells.Item(25, 6).Value = 3.360735064642597
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 13).Value = 0.4227369119627582
